$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 108661
$ws.Range("E2").Value = -343
$ws.Range("F2").Value = -343
$ws.Range("G2").Value = -2575
$ws.Range("H2").Value = -3207
$ws.Range("I2").Value = -3469
$ws.Range("J2").Value = 262
$ws.Range("K2").Value = 150728
$ws.Range("L2").Value = 78549
$ws.Range("M2").Value = 72179
$ws.Range("N2").Value = 59069
$ws.Range("O2").Value = 13110
$ws.Range("P2").Value = 4735
$ws.Range("Q2").Value = 5004
$ws.Range("R2").Value = -17577
$ws.Range("S2").Value = 12307
$ws.Range("T2").Value = 12960
$ws.Range("U2").Value = -7957
$ws.Range("V2").Value = 58788
$ws.Range("W2").Value = -0.32
$ws.Range("X2").Value = -2.95
$ws.Range("Y2").Value = -5.56
$ws.Range("Z2").Value = -2.31
$ws.Range("AA2").Value = 108.82
$ws.Range("AB2").Value = 1208.72
$ws.Range("AC2").Value = -3663
$ws.Range("AD2").Value = -10.92
$ws.Range("AE2").Value = 62391
$ws.Range("AF2").Value = 0.64
$ws.Range("AG2").Value = 1200
$ws.Range("AH2").Value = 3
$ws.Range("AI2").Value = -32.77
$ws.Range("AJ2").Value = 92915378
$ws.Range("D3").Value = 121795
$ws.Range("E3").Value = 16043
$ws.Range("F3").Value = 15818
$ws.Range("G3").Value = 13091
$ws.Range("H3").Value = 5093
$ws.Range("I3").Value = 4922
$ws.Range("J3").Value = 171
$ws.Range("K3").Value = 190230
$ws.Range("L3").Value = 109252
$ws.Range("M3").Value = 80979
$ws.Range("N3").Value = 63209
$ws.Range("O3").Value = 17770
$ws.Range("P3").Value = 4735
$ws.Range("Q3").Value = 6530
$ws.Range("R3").Value = -29441
$ws.Range("S3").Value = 22044
$ws.Range("T3").Value = 17766
$ws.Range("U3").Value = -11236
$ws.Range("V3").Value = 83717
$ws.Range("W3").Value = 13.17
$ws.Range("X3").Value = 4.18
$ws.Range("Y3").Value = 8.050000000000001
$ws.Range("Z3").Value = 2.99
$ws.Range("AA3").Value = 134.91
$ws.Range("AB3").Value = 1895.01
$ws.Range("AC3").Value = 5198
$ws.Range("AD3").Value = 9.75
$ws.Range("AE3").Value = 66764
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 2.96
$ws.Range("AI3").Value = 28.87
$ws.Range("AJ3").Value = 92915378
$ws.Range("D4").Value = 134624
$ws.Range("E4").Value = 17542
$ws.Range("F4").Value = 17542
$ws.Range("G4").Value = 16291
$ws.Range("H4").Value = 9233
$ws.Range("I4").Value = 8057
$ws.Range("J4").Value = 1176
$ws.Range("K4").Value = 203412
$ws.Range("L4").Value = 114209
$ws.Range("M4").Value = 89203
$ws.Range("N4").Value = 69642
$ws.Range("O4").Value = 19561
$ws.Range("P4").Value = 4735
$ws.Range("Q4").Value = 8517
$ws.Range("R4").Value = -5236
$ws.Range("S4").Value = -3363
$ws.Range("T4").Value = 9231
$ws.Range("U4").Value = -714
$ws.Range("V4").Value = 85234
$ws.Range("W4").Value = 13.03
$ws.Range("X4").Value = 6.86
$ws.Range("Y4").Value = 12.13
$ws.Range("Z4").Value = 4.69
$ws.Range("AA4").Value = 128.03
$ws.Range("AB4").Value = 2031.22
$ws.Range("AC4").Value = 8508
$ws.Range("AD4").Value = 6.36
$ws.Range("AE4").Value = 73559
$ws.Range("AF4").Value = 0.74
$ws.Range("AG4").Value = 1600
$ws.Range("AH4").Value = 2.96
$ws.Range("AI4").Value = 18.81
$ws.Range("AJ4").Value = 92915378
$ws.Range("D5").Value = 155801
$ws.Range("E5").Value = 20408
$ws.Range("F5").Value = 20408
$ws.Range("G5").Value = 18740
$ws.Range("H5").Value = 10919
$ws.Range("I5").Value = 9862
$ws.Range("J5").Value = 1057
$ws.Range("K5").Value = 220037
$ws.Range("L5").Value = 123411
$ws.Range("M5").Value = 96626
$ws.Range("N5").Value = 77302
$ws.Range("O5").Value = 19324
$ws.Range("P5").Value = 4735
$ws.Range("Q5").Value = 13065
$ws.Range("R5").Value = -14182
$ws.Range("S5").Value = 3345
$ws.Range("T5").Value = 11180
$ws.Range("U5").Value = 1885
$ws.Range("V5").Value = 90669
$ws.Range("W5").Value = 13.1
$ws.Range("X5").Value = 7.01
$ws.Range("Y5").Value = 13.42
$ws.Range("Z5").Value = 5.16
$ws.Range("AA5").Value = 127.72
$ws.Range("AB5").Value = 2190.74
$ws.Range("AC5").Value = 10414
$ws.Range("AD5").Value = 5.97
$ws.Range("AE5").Value = 81650
$ws.Range("AF5").Value = 0.76
$ws.Range("AG5").Value = 1800
$ws.Range("AH5").Value = 2.89
$ws.Range("AI5").Value = 17.29
$ws.Range("AJ5").Value = 92915378
$ws.Range("D6").Value = 177444
$ws.Range("E6").Value = 22098
$ws.Range("F6").Value = 22098
$ws.Range("G6").Value = 19397
$ws.Range("H6").Value = 10305
$ws.Range("I6").Value = 9032
$ws.Range("K6").Value = 226310
$ws.Range("L6").Value = 121060
$ws.Range("M6").Value = 105250
$ws.Range("N6").Value = 84904
$ws.Range("P6").Value = 4735
$ws.Range("Q6").Value = 12692
$ws.Range("R6").Value = -8359
$ws.Range("S6").Value = -5150
$ws.Range("T6").Value = 7313
$ws.Range("U6").Value = 5379
$ws.Range("V6").Value = 90110
$ws.Range("W6").Value = 12.45
$ws.Range("X6").Value = 5.81
$ws.Range("Y6").Value = 11.14
$ws.Range("Z6").Value = 4.62
$ws.Range("AA6").Value = 115.02
$ws.Range("AB6").Value = 2342.86
$ws.Range("AC6").Value = 9537
$ws.Range("AD6").Value = 5.41
$ws.Range("AE6").Value = 89680
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6").Value = 1900
$ws.Range("AH6").Value = 3.68
$ws.Range("AI6").Value = 19.93
$ws.Range("AJ6").Value = 92915378
$ws.Range("D7").Value = 179898
$ws.Range("E7").Value = 20914
$ws.Range("G7").Value = 19224
$ws.Range("H7").Value = 8714
$ws.Range("I7").Value = 7179
$ws.Range("K7").Value = 248024
$ws.Range("L7").Value = 136190
$ws.Range("M7").Value = 111834
$ws.Range("N7").Value = 90360
$ws.Range("P7").Value = 4739
$ws.Range("Q7").Value = 19581
$ws.Range("R7").Value = -6696
$ws.Range("S7").Value = -5464
$ws.Range("T7").Value = 7442
$ws.Range("U7").Value = 13944
$ws.Range("W7").Value = 11.63
$ws.Range("X7").Value = 4.84
$ws.Range("Y7").Value = 8.19
$ws.Range("Z7").Value = 3.67
$ws.Range("AA7").Value = 121.78
$ws.Range("AC7").Value = 7581
$ws.Range("AD7").Value = 6.03
$ws.Range("AE7").Value = 95442
$ws.Range("AF7").Value = 0.48
$ws.Range("AG7").Value = 1960
$ws.Range("AH7").Value = 4.29
$ws.Range("AI7").Value = 25.37
$ws.Range("D8").Value = 191207
$ws.Range("E8").Value = 23625
$ws.Range("G8").Value = 21630
$ws.Range("H8").Value = 10492
$ws.Range("I8").Value = 9787
$ws.Range("K8").Value = 257317
$ws.Range("L8").Value = 137610
$ws.Range("M8").Value = 119709
$ws.Range("N8").Value = 96848
$ws.Range("P8").Value = 4739
$ws.Range("Q8").Value = 11703
$ws.Range("R8").Value = -10077
$ws.Range("S8").Value = -2528
$ws.Range("T8").Value = 7676
$ws.Range("U8").Value = 2411
$ws.Range("W8").Value = 12.36
$ws.Range("X8").Value = 5.49
$ws.Range("Y8").Value = 10.46
$ws.Range("Z8").Value = 4.15
$ws.Range("AA8").Value = 114.95
$ws.Range("AC8").Value = 10334
$ws.Range("AD8").Value = 4.42
$ws.Range("AE8").Value = 102295
$ws.Range("AF8").Value = 0.45
$ws.Range("AG8").Value = 2070
$ws.Range("AH8").Value = 4.53
$ws.Range("AI8").Value = 19.65
$ws.Range("D9").Value = 195843
$ws.Range("E9").Value = 24270
$ws.Range("G9").Value = 21751
$ws.Range("H9").Value = 10470
$ws.Range("I9").Value = 9933
$ws.Range("K9").Value = 267400
$ws.Range("L9").Value = 140144
$ws.Range("M9").Value = 127253
$ws.Range("N9").Value = 102718
$ws.Range("P9").Value = 4739
$ws.Range("Q9").Value = 11599
$ws.Range("R9").Value = -9770
$ws.Range("S9").Value = -1028
$ws.Range("T9").Value = 7282
$ws.Range("U9").Value = 3077
$ws.Range("W9").Value = 12.39
$ws.Range("X9").Value = 5.35
$ws.Range("Y9").Value = 9.949999999999999
$ws.Range("Z9").Value = 3.99
$ws.Range("AA9").Value = 110.13
$ws.Range("AC9").Value = 10489
$ws.Range("AD9").Value = 4.36
$ws.Range("AE9").Value = 108495
$ws.Range("AF9").Value = 0.42
$ws.Range("AG9").Value = 2120
$ws.Range("AH9").Value = 4.64
$ws.Range("AI9").Value = 19.83
